$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (sheet row 1) gains a new "Password" column ---
$ws.Cells.Item(1,7).Value = "Password"

# --- "trashboatsr" row (sheet row 2): token rotated, new "akhil" password added ---
$ws.Cells.Item(2,5).Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiI0ODI5OGVhMC0yNDBhLTExZWUtOWMwNC1iMzcyMDk2MTViOGIiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiY2NiYTNlNWMtMGEyNy00MmIwLTgzNDUtNmE5MzQyNjFhMWEzIiwiaWF0IjoxNzA0MjQ3OTY5fQ.MQiQcDhDYkB2LxltxmzZZgIPhfxhnB6gFiVbMwAigCs"
$ws.Cells.Item(2,7).Value = "akhil"

# --- "RichDogeyBoy" row (sheet row 3) is removed; its data lives on in row 0 already ---
$ws.Rows.Item(3).Delete()
